$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45251
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = "$/bandeja 10 kilos"
$ws.Range("R2").Value = "Provincia de Quillota"
$ws.Range("S2").Value = 1500
$ws.Range("T2").Value = 10

# Row 3
$ws.Range("D3").Value = 45251
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 40
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("R3").Value = "Provincia de Quillota"
$ws.Range("S3").Value = 1200
$ws.Range("T3").Value = 10

# Row 4
$ws.Range("D4").Value = 45251
$ws.Range("L4").Value = "Tercera"
$ws.Range("M4").Value = 35
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 10000
$ws.Range("Q4").Value = "$/bandeja 10 kilos"
$ws.Range("R4").Value = "Provincia de Quillota"
$ws.Range("S4").Value = 1000
$ws.Range("T4").Value = 10

# Row 5
$ws.Range("D5").Value = 44911
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 220
$ws.Range("N5").Value = 5000
$ws.Range("O5").Value = 5000
$ws.Range("P5").Value = 5000
$ws.Range("Q5").Value = "$/bandeja 5 kilos"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1000
$ws.Range("T5").Value = 5

# Row 6
$ws.Range("D6").Value = 44911
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 4000
$ws.Range("O6").Value = 4000
$ws.Range("P6").Value = 4000
$ws.Range("Q6").Value = "$/bandeja 5 kilos"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 800
$ws.Range("T6").Value = 5

# Row 7
$ws.Range("D7").Value = 44915
$ws.Range("L7").Value = "Especial"
$ws.Range("M7").Value = 150
$ws.Range("N7").Value = 6000
$ws.Range("O7").Value = 6000
$ws.Range("P7").Value = 6000
$ws.Range("Q7").Value = "$/bandeja 5 kilos"
$ws.Range("R7").Value = "Provincia de Quillota"
$ws.Range("S7").Value = 1200
$ws.Range("T7").Value = 5

# Row 8
$ws.Range("D8").Value = 44915
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 5000
$ws.Range("O8").Value = 5000
$ws.Range("P8").Value = 5000
$ws.Range("Q8").Value = "$/bandeja 5 kilos"
$ws.Range("R8").Value = "Provincia de Quillota"
$ws.Range("S8").Value = 1000
$ws.Range("T8").Value = 5

